# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.960.69'
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').Value = '1.634.84'
$ws.Range('E3').Value = '  -0.54%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''214.21'
$ws.Range('E6').Value = '  -0.67%  '
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('E9').Value = '  -2.43%  '
$ws.Range('D10').Value = '''18.51'
$ws.Range('E10').Value = '  -5.85%  '
$ws.Range('E11').Value = '  -0.74%  '
$ws.Range('D12').Value = '1.862.62'
$ws.Range('E12').Value = '  -0.53%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.19'
$ws.Range('E13').Value = '  -2.19%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.588.35'
$ws.Range('E14').Value = '  -4.24%  '
$ws.Range('D15').Value = '''0.530'
$ws.Range('E15').Value = '  -2.64%  '
$ws.Range('D16').Value = '25.975.39'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('D17').Value = '0.0₃0743'
$ws.Range('E17').Value = '  -2.89%  '
$ws.Range('D18').Value = '''61.76'
$ws.Range('E18').Value = '  -2.32%  '
$ws.Range('E19').Value = '  -0.17%  '
$ws.Range('D20').Value = '''190.91'
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('E21').Value = '  -2.44%  '
$ws.Range('D22').Value = '''9.71'
$ws.Range('E22').Value = '  -2.40%  '
$ws.Range('E23').Value = '  -2.24%  '
$ws.Range('E24').Value = '  +1.24%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').Value = '''143.62'
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D26').Value = '''1.78'
$ws.Range('E26').Value = '  -0.87%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('D28').Value = '''6.84'
$ws.Range('E28').Value = '  -0.82%  '
$ws.Range('E29').Value = '  -1.78%  '
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').Value = '''0.0483'
$ws.Range('E31').Value = '  -3.09%  '
$ws.Range('D32').Value = '''3.15'
$ws.Range('E32').Value = '  -3.20%  '
$ws.Range('E33').Value = '  -4.32%  '
$ws.Range('E34').Value = '  -2.55%  '
$ws.Range('E35').Value = '  -2.75%  '
$ws.Range('D36').Value = '1.134.61'
$ws.Range('E36').Value = '  +0.26%  '
$ws.Range('D37').Value = '''0.865'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('E38').Value = '  -1.33%  '
$ws.Range('D39').Value = '''0.520'
$ws.Range('E39').Value = '  -3.74%  '
$ws.Range('E40').Value = '  -1.61%  '
$ws.Range('E41').Value = '  -1.09%  '
$ws.Range('E42').Value = '  -2.48%  '
$ws.Range('E43').Value = '  -0.54%  '
$ws.Range('E44').Value = '  -4.88%  '
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').Value = '''55.19'
$ws.Range('E46').Value = '  -2.34%  '
$ws.Range('E47').Value = '  -0.25%  '
$ws.Range('E48').Value = '  +2.44%  '
$ws.Range('E49').Value = '  -0.54%  '
$ws.Range('D50').Value = '''7.53'
$ws.Range('E50').Value = '  -2.57%  '
$ws.Range('E51').Value = '  -0.03%  '
